# "creator statement and project reflection" — apply the author's edits
# using Word COM-interop Find/Replace against the whole document story.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- Paragraph: "In my project, I tell a story ..." ---
Replace-Text `
    "In my project, I tell a story from my childhood golf summer camp tournament over the sounds of me practicing at my local driving range. I used a Zoom H5 Handy Recorder from the Hillman Library to capture my audio clips and a Logitech Blue Yeti USB microphone to record my narration. I layered appropriate sounds over different parts of the story, including practice swings, drives, and putts. My audio includes birds chirping, muffled music, and other people conversing at the range. I like how the background noise creates a fast-paced, succinct tone that merges well with the pace of my narration." `
    "In my project, I narrate a story from my childhood golf summer camp tournament over the sounds of me practicing at my local driving range. I used a Zoom H5 Handy Recorder from the Hillman Library to capture my audio clips and a Logitech Blue Yeti USB microphone to record my narration. I layered appropriate sounds over different parts of the story, including practice swings, drives, and putts. My audio includes birds chirping, muffled music, and other people conversing at the range. The background noise creates a fast-paced, succinct tone that merges well with the pace of my narration."

# --- Paragraph: "When I got home and reviewed my clips ..." ---
Replace-Text `
    "When I got home and reviewed my clips, it occurred to me that my brain filters out a lot of everyday background noise. My project suddenly became a lot more challenging. A lot of trucks and cars go unnoticed, as well as far-off conversations. The microphone I used was good, and it picked up everything. I knew my job editing would be easier the cleaner the audio I captured." `
    "When I got home and reviewed my clips, it occurred to me that my brain filters out a lot of everyday background noise. I knew my job editing would be easier if I could capture audio with minimal background noise. I ended up liking the way my narration sounded including the background noises at the range."

# --- Paragraph: "My next task was to combine my clips ..." ---
Replace-Text `
    "My next task was to combine my clips to tell a story. At this point, I only had sounds of me driving. I knew I would need more, maybe I could organize my clips in descending fashion, from driving, to chipping, to putting. This would wind the audience down from the loudest sounds of golf to the quiet soothing ones." `
    "My next task was to combine my clips to tell a story. At this point, I only had sounds of me driving. I knew I would need more, maybe I could organize my clips in descending order, from driving, to chipping, to putting. This would wind the audience down from the loudest sounds of golf to the more soothing ones, giving them a sense that the story would conclude with putting."

# --- Paragraph: "I decided to try to record other sounds ..." ---
Replace-Text `
    "I decided to try to record other sounds, such as making coffee and brushing my teeth, to give the audience more of a sense of daily routine. These clips were difficult to record and made my project more confusing, so I scrapped them. I found that boiling my project down to the essence of golf painted a clearer picture of time, place, and setting. I have learned making errors is very important to learn how to make successes. I embraced some of the errors that occurred over the week, but some challenges were flat out annoying and frustrating." `
    "I decided to experiment with recording other sounds, such as making coffee and brushing my teeth, to give the audience more of a sense of daily routine. These clips were difficult to record and made my project more confusing, so I scrapped them. I found that editing my project down to the essence of golf painted a clearer picture of time, place, and setting. I have learned making errors is very important to learn how to make successes. I embraced some of the errors that occurred over the week, and some were flat out annoying and frustrating."

# --- Paragraph: "I was still hesitant about my golf vision ..." ---
Replace-Text `
    "I was still hesitant about my golf vision, so I went on a hike and recorded more sounds, and most of the clips were unusable. There were a lot of unwanted sounds, such as microphone rubbing and wind. If I attempted to record a hike again, I would bring a set of headphones to listen to make sure I was capturing what I intended." `
    "I was still hesitant about my golf vision, so I went on a hike and recorded more sounds. Most of the clips were unusable. There were a lot of unwanted sounds, such as accidentally rubbing the microphone and wind. If I attempted to record a hike again, I would bring a set of headphones to listen to make sure I was capturing what I intended."

# --- Sentence inside the "As the deadline neared ..." paragraph ---
Replace-Text `
    "I did end up taking one sentence from the first take and replacing it with the exact sentence from the second take." `
    "I took one sentence from the first take and replaced it with the same sentence from the second take."
